$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()

# Enter the value into E8
$ws2.Range("E8").Value = "E8"

# Set up the selection before freezing panes
$ws2.Range("G14").Select()

# Freeze panes at E8 (xSplit=4, ySplit=7)
$ws2.Range("E8").Select()
$excel.ActiveWindow.FreezePanes = $true

# Restore the active cell selection in the bottom-right pane
$ws2.Range("G14").Select()
